$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-7 with recomputed NATMI values (new TPM input) and
# reordered Sending/Target cluster + Ligand/Receptor symbol strings.
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Wnt5a"
$ws.Cells.Item(2,3).Value = "Fzd6"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.068624
$ws.Cells.Item(2,8).Value = 0.205872
$ws.Cells.Item(2,9).Value = 0.01198115042951486
$ws.Cells.Item(2,10).Value = 0.01198115042951486
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 12.36292333333333
$ws.Cells.Item(2,14).Value = 37.08877
$ws.Cells.Item(2,15).Value = 0.918273862214392
$ws.Cells.Item(2,16).Value = 0.9303794466068031
$ws.Cells.Item(2,17).Value = 0.8483932508266667
$ws.Cells.Item(2,18).Value = 7.63553925744
$ws.Cells.Item(2,19).Value = 0.01100197727868223
$ws.Cells.Item(2,20).Value = 0.0111470161063249
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Wnt5a"
$ws.Cells.Item(3,3).Value = "Fzd6"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.068624
$ws.Cells.Item(3,8).Value = 0.205872
$ws.Cells.Item(3,9).Value = 0.01198115042951486
$ws.Cells.Item(3,10).Value = 0.01198115042951486
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.574769
$ws.Cells.Item(3,14).Value = 1.724307
$ws.Cells.Item(3,15).Value = 0.04269179184247177
$ws.Cells.Item(3,16).Value = 0.04325459680761149
$ws.Cells.Item(3,17).Value = 0.039442947856
$ws.Cells.Item(3,18).Value = 0.354986530704
$ws.Cells.Item(3,19).Value = 0.0005114967801701896
$ws.Cells.Item(3,20).Value = 0.0005182398311200066
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Wnt5a"
$ws.Cells.Item(4,3).Value = "Fzd6"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.068624
$ws.Cells.Item(4,8).Value = 0.205872
$ws.Cells.Item(4,9).Value = 0.01198115042951486
$ws.Cells.Item(4,10).Value = 0.01198115042951486
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 0.525528
$ws.Cells.Item(4,14).Value = 1.051056
$ws.Cells.Item(4,15).Value = 0.03903434594313629
$ws.Cells.Item(4,16).Value = 0.02636595658558534
$ws.Cells.Item(4,17).Value = 0.03606383347200001
$ws.Cells.Item(4,18).Value = 0.216383000832
$ws.Cells.Item(4,19).Value = 0.000467676370662439
$ws.Cells.Item(4,20).Value = 0.000315894492069956
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Wnt5a"
$ws.Cells.Item(5,3).Value = "Fzd6"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 5.659039666666668
$ws.Cells.Item(5,8).Value = 16.977119
$ws.Cells.Item(5,9).Value = 0.9880188495704851
$ws.Cells.Item(5,10).Value = 0.9880188495704851
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 12.36292333333333
$ws.Cells.Item(5,14).Value = 37.08877
$ws.Cells.Item(5,15).Value = 0.918273862214392
$ws.Cells.Item(5,16).Value = 0.9303794466068031
$ws.Cells.Item(5,17).Value = 69.96227353929223
$ws.Cells.Item(5,18).Value = 629.6604618536301
$ws.Cells.Item(5,19).Value = 0.9072718849357098
$ws.Cells.Item(5,20).Value = 0.9192324305004782
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Wnt5a"
$ws.Cells.Item(6,3).Value = "Fzd6"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 5.659039666666668
$ws.Cells.Item(6,8).Value = 16.977119
$ws.Cells.Item(6,9).Value = 0.9880188495704851
$ws.Cells.Item(6,10).Value = 0.9880188495704851
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.574769
$ws.Cells.Item(6,14).Value = 1.724307
$ws.Cells.Item(6,15).Value = 0.04269179184247177
$ws.Cells.Item(6,16).Value = 0.04325459680761149
$ws.Cells.Item(6,17).Value = 3.252640570170334
$ws.Cells.Item(6,18).Value = 29.273765131533
$ws.Cells.Item(6,19).Value = 0.04218029506230157
$ws.Cells.Item(6,20).Value = 0.04273635697649148
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Wnt5a"
$ws.Cells.Item(7,3).Value = "Fzd6"
$ws.Cells.Item(7,4).Value = "MuSCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 5.659039666666668
$ws.Cells.Item(7,8).Value = 16.977119
$ws.Cells.Item(7,9).Value = 0.9880188495704851
$ws.Cells.Item(7,10).Value = 0.9880188495704851
$ws.Cells.Item(7,11).Value = 2
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.525528
$ws.Cells.Item(7,14).Value = 1.051056
$ws.Cells.Item(7,15).Value = 0.03903434594313629
$ws.Cells.Item(7,16).Value = 0.02636595658558534
$ws.Cells.Item(7,17).Value = 2.973983797944
$ws.Cells.Item(7,18).Value = 17.843902787664
$ws.Cells.Item(7,19).Value = 0.03856666957247385
$ws.Cells.Item(7,20).Value = 0.02605006209351539

# The old rows 8-10 (ECs -> MuSCs / FAPs / ECs ligand-receptor combos using
# the previous cluster ordering) are no longer present in the recomputed
# output; remove them so the sheet matches the regenerated table.
$ws.Range("A8:A10").EntireRow.Delete()
